$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L3 value from 2020 to 2021
$ws.Range("L3").Value = 2021

# Add new column M data for row 3 and row 4, matching style of existing L column cells
$ws.Range("M3").Value = 2022
$ws.Range("M4").Value = 6.18

# Copy styles from L3/L4 to M3/M4 so number formatting / borders match
$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update selection to reflect final active cell M9
$ws.Range("M9").Select() | Out-Null
